$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 5's data)
$ws.Range("A2").Value = 111525233
$ws.Range("B2").Value = 77515
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("Q2").Value = 404540.9329893424
$ws.Range("R2").Value = 6706716.233959051

# Row 3 (was row 6's data)
$ws.Range("A3").Value = 111525238
$ws.Range("B3").Value = 77515
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("Q3").Value = 404495.4563026094
$ws.Range("R3").Value = 6706677.491168984

# Row 5 (was row 7's data)
$ws.Range("A5").Value = 111525226
$ws.Range("B5").Value = 77515
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("Q5").Value = 404616.9589749529
$ws.Range("R5").Value = 6706770.937089294

# Row 6 (was row 3's data)
$ws.Range("A6").Value = 111525224
$ws.Range("B6").Value = 77268
$ws.Range("E6").Value = 228912
$ws.Range("F6").Value = "Mörk kolflarnlav"
$ws.Range("G6").Value = "Carbonicola myrmecina"
$ws.Range("H6").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q6").Value = 404619.9854206198
$ws.Range("R6").Value = 6706773.322858612

# Row 7 (was row 2's data)
$ws.Range("A7").Value = 111525223
$ws.Range("B7").Value = 78107
$ws.Range("E7").Value = 6453
$ws.Range("F7").Value = "Vedskivlav"
$ws.Range("G7").Value = "Hertelidea botryosa"
$ws.Range("H7").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q7").Value = 404637.0659126193
$ws.Range("R7").Value = 6706784.214121711
